# Update CDA Logical model for ST.r2b
# - rename the "Include" sheet tab
# - bump Version / Date metadata
# - insert a "Jurisdiction" metadata row (empty value) before "Description"

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)
$wsInclude = $wb.Worksheets.Item(2)

# 1. Rename the include sheet
$wsInclude.Name = "Include #0"

# 2. Bump the Version and Date values on the Metadata sheet
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 3. Insert a new "Jurisdiction" row right after "Contact" (row 10) and
#    before "Description" (row 11), pushing everything else down by one.
$wsMeta.Rows.Item(11).Insert()

$wsMeta.Range("A11").Value = "Jurisdiction"
# Leading apostrophe forces a literal (empty) text value instead of clearing the cell
$wsMeta.Range("B11").Value = "'"

# Match the body-row formatting (border/wrap) used by the rest of the table
# (applied after the values so the format copy wins over the quote-prefix flag)
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
